$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename row-3 description text (PLP details -> Product List Page Scenarios, trailing space kept)
$ws.Range("B3").Value = "Product List Page Scenarios "

# New row 4 - set description (B) before TSID (A) so the shared-string table order matches
$ws.Range("B4").Value = "Product Details Page Scenarios"
$ws.Range("A4").Value = "ProductDetailsPage"
$ws.Range("C4").Value = "Y"

# Copy the formatting from row 3 down onto the new row 4 cells
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)

# Remove the word-wrap on B2 (creates the new cellXfs entry)
$ws.Range("B2").WrapText = $false

# Move the active selection to the newly added row
[void]$ws.Range("A4").Select()
